# Final Changes 31 Dec
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Sheet1: change A2 from "excelTest" to "testSignUpFirstPage"
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "testSignUpFirstPage"

# ---------------------------------------------------------------------------
# Add Sheet2 ("testLoginPage") right after Sheet1
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "TestCaseName"
$ws2.Range("B1").Value = "VendorFirstName"
$ws2.Range("C1").Value = "VendorLastName"

$ws2.Range("A2").Value = "testLoginPage"
$ws2.Range("B2").Value = "apandhe5@xpanxion.com"
$ws2.Range("C2").Value = "Test#123"

$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:apandhe5@xpanxion.com")
$ws2.Range("B2").Style = "Hyperlink"

$ws2.Range("A1").ColumnWidth = 14.28515625
$ws2.Range("B1").ColumnWidth = 24.42578125
$ws2.Range("C1").ColumnWidth = 16.42578125

$ws2.Range("A1:C2").Select()

# ---------------------------------------------------------------------------
# Add Sheet3 ("testAddApplicantPage") right after Sheet2
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

$ws3.Range("A1").Value = "TestCaseName"
$ws3.Range("B1").Value = "VendorFirstName"
$ws3.Range("C1").Value = "VendorLastName"
$ws3.Range("D1").Value = "ApplicantFirstName"
$ws3.Range("E1").Value = "ApplicantLastName"
$ws3.Range("F1").Value = "ApplicantAddressLine1"
$ws3.Range("G1").Value = "ApplicantCountry"
$ws3.Range("H1").Value = "ApplicantCity"
$ws3.Range("I1").Value = "ApplicantState"
$ws3.Range("J1").Value = "ApplicantPostalCode"
$ws3.Range("K1").Value = "ApplicantEmailAddress"
$ws3.Range("L1").Value = "ApplicantDOB"
$ws3.Range("M1").Value = "ApplicantConfirmDOB"
$ws3.Range("N1").Value = "ApplicantSSN"
$ws3.Range("O1").Value = "ApplicantConfirmSSN"

$ws3.Range("A2").Value = "testAddApplicantPage"
$ws3.Range("B2").Value = "apandhe5@xpanxion.com"
$ws3.Range("C2").Value = "Test#123"
$ws3.Range("D2").Value = "Max"
$ws3.Range("E2").Value = "Min"
$ws3.Range("F2").Value = "SrNo"
$ws3.Range("G2").Value = "Sr No"
$ws3.Range("H2").Value = "NY"
$ws3.Range("I2").Value = "NY"
$ws3.Range("J2").Value = 10001
$ws3.Range("K2").Value = "aps@gps.com"
$ws3.Range("L2").Value = 32874
$ws3.Range("L2").NumberFormat = "mm-dd-yy"
$ws3.Range("M2").Value = 32874
$ws3.Range("M2").NumberFormat = "mm-dd-yy"
$ws3.Range("N2").Value = 456999990
$ws3.Range("O2").Value = 456999990

$ws3.Hyperlinks.Add($ws3.Range("B2"), "mailto:apandhe5@xpanxion.com")
$ws3.Range("B2").Style = "Hyperlink"
$ws3.Hyperlinks.Add($ws3.Range("K2"), "mailto:aps@gps.com")
$ws3.Range("K2").Style = "Hyperlink"

$ws3.Range("A1").ColumnWidth = 21
$ws3.Range("B1").ColumnWidth = 24.42578125
$ws3.Range("C1").ColumnWidth = 16.42578125
$ws3.Range("D1").ColumnWidth = 18.85546875
$ws3.Range("E1").ColumnWidth = 18.42578125
$ws3.Range("F1").ColumnWidth = 21.7109375
$ws3.Range("G1").ColumnWidth = 16.5703125
$ws3.Range("H1").ColumnWidth = 12.85546875
$ws3.Range("I1").ColumnWidth = 14.140625
$ws3.Range("J1").ColumnWidth = 19.7109375
$ws3.Range("K1").ColumnWidth = 21.85546875
$ws3.Range("L1").ColumnWidth = 13.42578125
$ws3.Range("M1").ColumnWidth = 20.7109375
$ws3.Range("N1").ColumnWidth = 12.85546875
$ws3.Range("O1").ColumnWidth = 20.28515625

# ---------------------------------------------------------------------------
# Selections / active sheet: Sheet3 ends up active with B14 selected;
# Sheet1's prior tabSelected/topLeftCell goes away, selection becomes B17.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B17").Select()

$ws3.Activate()
$ws3.Range("B14").Select()
